$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Giovanni"
$ws.Range("B10").Value = "Interno"
$ws.Range("C10").Value = 32

$ws.Range("D9").Copy($ws.Range("D10"))
$ws.Range("D10").Value = 43499

$ws.Range("C11").Select()
